# Update "想去人数" (want-to-go count) figures for several rows across
# the "展览" and "全部类型" worksheets, matching a refreshed data scrape.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 128
$ws1.Range("F10").Value = 15
$ws1.Range("F16").Value = 4385
$ws1.Range("F20").Value = 221
$ws1.Range("F25").Value = 3180
$ws1.Range("F37").Value = 5587
$ws1.Range("F38").Value = 857
$ws1.Range("F42").Value = 50
$ws1.Range("F43").Value = 1126
$ws1.Range("F44").Value = 502
$ws1.Range("F46").Value = 2007

# --- 演出 (sheet 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 64
$ws2.Range("F16").Value = 126

# --- 全部类型 (sheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 128
$ws4.Range("F12").Value = 64
$ws4.Range("F17").Value = 4385
$ws4.Range("F21").Value = 221
$ws4.Range("F23").Value = 3180
$ws4.Range("F34").Value = 126
$ws4.Range("F36").Value = 5587
$ws4.Range("F38").Value = 857
$ws4.Range("F44").Value = 50
$ws4.Range("F45").Value = 1126
$ws4.Range("F46").Value = 502
$ws4.Range("F47").Value = 2007
